# chore: update Sheets via scheduled runner
# Refreshes market-board-derived price/profit figures (columns H-N) for a
# handful of leve rows across each job sheet, leaving all other data
# (item names, leve metadata, etc.) untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39950
$ws.Range("J3").Value = 39950
$ws.Range("L3").Value = 39950
$ws.Range("N3").Value = -40178

$ws.Range("H93").Value = 50601
$ws.Range("J93").Value = 50601
$ws.Range("L93").Value = 50601
$ws.Range("N93").Value = -55593

$ws.Range("H102").Value = 39950
$ws.Range("J102").Value = 39950
$ws.Range("L102").Value = 39950
$ws.Range("N102").Value = -46440

$ws.Range("H113").Value = 2550.7144
$ws.Range("I113").Value = 1988.75
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 1988.75
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = 1265.25
$ws.Range("N113").Value = -9808

$ws.Range("H137").Value = 1724.6052
$ws.Range("I137").Value = 1853.9445
$ws.Range("J137").Value = 1608.2
$ws.Range("K137").Value = 5561.833500000001
$ws.Range("L137").Value = 4824.6
$ws.Range("M137").Value = -3011.833500000001
$ws.Range("N137").Value = -9924.6

$ws.Range("H138").Value = 1959.09
$ws.Range("I138").Value = 1119.9565
$ws.Range("J138").Value = 2209.7402
$ws.Range("K138").Value = 3359.8695
$ws.Range("L138").Value = 6629.220600000001
$ws.Range("M138").Value = 1780.1305
$ws.Range("N138").Value = -16909.2206

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 826.4545000000001
$ws.Range("I4").Value = 809.1
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 809.1
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -693.1
$ws.Range("N4").Value = -1232

$ws.Range("H45").Value = 2499.8333
$ws.Range("I45").Value = 2399.8
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 2399.8
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -2022.8
$ws.Range("N45").Value = -3754

$ws.Range("H92").Value = 500000
$ws.Range("J92").Value = 500000
$ws.Range("L92").Value = 500000
$ws.Range("N92").Value = -504992

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null

$ws.Range("H132").Value = 1673.1177
$ws.Range("I132").Value = 1465.25
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4395.75
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1865.75
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2072.2222
$ws.Range("I99").Value = 1100
$ws.Range("J99").Value = 2193.75
$ws.Range("K99").Value = 1100
$ws.Range("L99").Value = 2193.75
$ws.Range("M99").Value = 398
$ws.Range("N99").Value = -5189.75

$ws.Range("H134").Value = 3487.4
$ws.Range("I134").Value = 3450.7856
$ws.Range("K134").Value = 10352.3568
$ws.Range("M134").Value = -7817.356800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 269821
$ws.Range("J28").Value = 269821
$ws.Range("L28").Value = 269821
$ws.Range("N28").Value = -270311

$ws.Range("H99").Value = 2038.9584
$ws.Range("I99").Value = 2273.9443
$ws.Range("K99").Value = 2273.9443
$ws.Range("M99").Value = -775.9443000000001

$ws.Range("H126").Value = 2038.9584
$ws.Range("I126").Value = 2273.9443
$ws.Range("K126").Value = 6821.8329
$ws.Range("M126").Value = -4351.8329

$ws.Range("H132").Value = 2499.3333
$ws.Range("I132").Value = 2226.1667
$ws.Range("J132").Value = 3318.8333
$ws.Range("K132").Value = 6678.500100000001
$ws.Range("L132").Value = 9956.499899999999
$ws.Range("M132").Value = -4148.500100000001
$ws.Range("N132").Value = -15016.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 967.8570999999999
$ws.Range("I92").Value = 999.5
$ws.Range("J92").Value = 925.6667
$ws.Range("K92").Value = 2998.5
$ws.Range("L92").Value = 2777.0001
$ws.Range("M92").Value = -1750.5
$ws.Range("N92").Value = -5273.0001

$ws.Range("H97").Value = 580
$ws.Range("I97").Value = 450
$ws.Range("J97").Value = 666.6667
$ws.Range("K97").Value = 1350
$ws.Range("L97").Value = 2000.0001
$ws.Range("M97").Value = -854
$ws.Range("N97").Value = -2992.0001

$ws.Range("H123").Value = 5550
$ws.Range("J123").Value = 10000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -34900

$ws.Range("H131").Value = 877.9693600000001
$ws.Range("J131").Value = 889.8
$ws.Range("L131").Value = 2669.4
$ws.Range("N131").Value = -12749.4

$ws.Range("H132").Value = 1697.5555
$ws.Range("I132").Value = 1015.6667
$ws.Range("J132").Value = 2379.4443
$ws.Range("K132").Value = 9141.0003
$ws.Range("L132").Value = 21414.9987
$ws.Range("M132").Value = -6611.0003
$ws.Range("N132").Value = -26474.9987

$ws.Range("H137").Value = 25642706
$ws.Range("I137").Value = 889.55554
$ws.Range("K137").Value = 2668.66662
$ws.Range("M137").Value = 2431.33338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 36500
$ws.Range("I20").Value = 36500
$ws.Range("K20").Value = 36500
$ws.Range("M20").Value = -36255

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null

$ws.Range("H69").Value = 179857.14
$ws.Range("J69").Value = 179857.14
$ws.Range("L69").Value = 179857.14
$ws.Range("N69").Value = -181355.14

$ws.Range("H70").Value = 6683.2607
$ws.Range("I70").Value = 6213.375
$ws.Range("J70").Value = 6933.8667
$ws.Range("K70").Value = 6213.375
$ws.Range("L70").Value = 6933.8667
$ws.Range("M70").Value = -5943.375
$ws.Range("N70").Value = -7473.8667

$ws.Range("H72").Value = 179857.14
$ws.Range("J72").Value = 179857.14
$ws.Range("L72").Value = 539571.42
$ws.Range("N72").Value = -547059.42

$ws.Range("H73").Value = 6683.2607
$ws.Range("I73").Value = 6213.375
$ws.Range("J73").Value = 6933.8667
$ws.Range("K73").Value = 6213.375
$ws.Range("L73").Value = 6933.8667
$ws.Range("M73").Value = -5277.375
$ws.Range("N73").Value = -8805.866699999999

$ws.Range("H126").Value = 2303.2273
$ws.Range("I126").Value = 2153.75
$ws.Range("J126").Value = 2701.8333
$ws.Range("K126").Value = 6461.25
$ws.Range("L126").Value = 8105.499899999999
$ws.Range("M126").Value = -3991.25
$ws.Range("N126").Value = -13045.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3499.6
$ws.Range("I7").Value = 2750.1667
$ws.Range("J7").Value = 4623.75
$ws.Range("K7").Value = 2750.1667
$ws.Range("L7").Value = 4623.75
$ws.Range("M7").Value = -2638.1667
$ws.Range("N7").Value = -4847.75

$ws.Range("H36").Value = 44800
$ws.Range("J36").Value = 44800
$ws.Range("L36").Value = 44800
$ws.Range("N36").Value = -45924

$ws.Range("H40").Value = 6984
$ws.Range("I40").Value = 9251
$ws.Range("J40").Value = 2450
$ws.Range("K40").Value = 9251
$ws.Range("L40").Value = 2450
$ws.Range("M40").Value = -9115
$ws.Range("N40").Value = -2722

$ws.Range("H46").Value = 1846.1538
$ws.Range("I46").Value = 1533.3334
$ws.Range("J46").Value = 1940
$ws.Range("K46").Value = 1533.3334
$ws.Range("L46").Value = 1940
$ws.Range("M46").Value = -1345.3334
$ws.Range("N46").Value = -2316

$ws.Range("H122").Value = 17863200
$ws.Range("I122").Value = 25006700
$ws.Range("J122").Value = 4450
$ws.Range("K122").Value = 75020100
$ws.Range("L122").Value = 13350
$ws.Range("M122").Value = -75017650
$ws.Range("N122").Value = -18250

$ws.Range("H126").Value = 3499.6
$ws.Range("I126").Value = 2750.1667
$ws.Range("J126").Value = 4623.75
$ws.Range("K126").Value = 8250.500100000001
$ws.Range("L126").Value = 13871.25
$ws.Range("M126").Value = -5780.500100000001
$ws.Range("N126").Value = -18811.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 107500

$ws.Range("H73").Value = 107500

$ws.Range("H120").Value = 31566.334
$ws.Range("J120").Value = 31566.334
$ws.Range("L120").Value = 31566.334
$ws.Range("N120").Value = -41242.334

$ws.Range("H132").Value = 2360.7222
$ws.Range("I132").Value = 2524.5625
$ws.Range("K132").Value = 7573.6875
$ws.Range("M132").Value = -5043.6875
